$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new blank paragraph right before the "nb-anacondacloud"
#    paragraph (i.e. right after "dynd==0.7.3.dev1").
# ------------------------------------------------------------------
$nbAnaconda = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "nb-anacondacloud*") {
        $nbAnaconda = $p
        break
    }
}
$nbAnaconda.Range.InsertParagraphBefore()

# Re-find the paragraph (collection shifted after the insert).
$nbAnaconda = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "nb-anacondacloud*") {
        $nbAnaconda = $p
        break
    }
}

# ------------------------------------------------------------------
# 2) Drop the stray <w:lastRenderedPageBreak/> hint that sits in front
#    of the "nb-anacondacloud" run. Rewriting the run's text (to a
#    throw-away value and back) forces the run to be re-emitted
#    without the rendering hint while keeping the spell-check markup
#    untouched.
# ------------------------------------------------------------------
$rng = $nbAnaconda.Range
$firstWordEnd = $rng.Start + 16   # length of "nb-anacondacloud"
$firstWordRange = $d.Range($rng.Start, $firstWordEnd)
$firstWordRange.Text = "TEMP_PLACEHOLDER"
$firstWordRange = $d.Range($rng.Start, $rng.Start + 17)
$firstWordRange.Text = "nb-anacondacloud"

# ------------------------------------------------------------------
# 3) Turn the old single-run "nb-conda-kernels==1.0.3" paragraph into
#    the properly spell-checked run sequence, and append a brand new
#    "nbpresent==3.0.2" paragraph (carrying the _GoBack bookmark that
#    used to live on the kernels paragraph).
# ------------------------------------------------------------------
$kernelsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "nb-conda-kernels*") {
        $kernelsPara = $p
        break
    }
}

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$xml = "<w:p $ns>" + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>nb</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>-</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>conda</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>-kernels==1.0.3</w:t></w:r>' + `
    '</w:p>' + `
    "<w:p $ns>" + `
    '<w:r><w:t>nbpresent==3.0.2</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    '</w:p>'

$kernelsPara.Range.InsertXML($xml)
